$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-text numbers (some using dotted
# thousands separators, e.g. "36.456.27"), not real numeric values. For the
# new prices that parse as a plain number, force the cell to Text format
# first so Excel does not silently convert the assigned string into a
# floating-point number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.456.27'
$ws.Range("E2").Value = '  -2.53%  '
$ws.Range("D3").Value = '1.985.10'
$ws.Range("E3").Value = '  -3.25%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '244.29'
$ws.Range("E5").Value = '  -3.05%  '
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("D7").Value = '59.02'
$ws.Range("E7").Value = '  -9.75%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.376'
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").Value = '57.70'
$ws.Range("E10").Value = '  -3.49%  '
$ws.Range("D11").Value = '0.0820'
$ws.Range("E11").Value = '  +6.74%  '
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").Value = '23.60'
$ws.Range("E13").Value = '  +6.29%  '
$ws.Range("D14").Value = '0.865'
$ws.Range("E14").Value = '  -5.31%  '
$ws.Range("D15").Value = '14.01'
$ws.Range("E15").Value = '  -5.30%  '
$ws.Range("D16").Value = '2.274.38'
$ws.Range("D17").Value = '5.47'
$ws.Range("E17").Value = '  -1.93%  '
$ws.Range("D18").Value = '1.981.66'
$ws.Range("E18").Value = '  -3.45%  '
$ws.Range("D19").Value = '36.345.06'
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").Value = '70.65'
$ws.Range("E20").Value = '  -4.17%  '
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("D22").Value = '5.34'
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").Value = '233.78'
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '2.59'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  -3.97%  '
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("D28").Value = '161.68'
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("D29").Value = '19.85'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  +11.20%  '
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("E33").Value = '  -7.10%  '
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").Value = '4.43'
$ws.Range("E35").Value = '  -5.68%  '
$ws.Range("D36").Value = '6.31'
$ws.Range("E36").Value = '  +4.72%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -7.22%  '
$ws.Range("E39").Value = '  -3.81%  '
$ws.Range("D40").Value = '3.05'
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("E42").Value = '  -7.79%  '
$ws.Range("D43").Value = '2.91'
$ws.Range("E43").Value = '  -3.74%  '
$ws.Range("D44").Value = '0.0214'
$ws.Range("E44").Value = '  -1.96%  '
$ws.Range("E45").Value = '  -4.69%  '
$ws.Range("D46").Value = '92.72'
$ws.Range("E46").Value = '  -4.13%  '
$ws.Range("D47").Value = '16.22'
$ws.Range("E47").Value = '  -5.22%  '
$ws.Range("D48").Value = '1.385.03'
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("D49").Value = '7.53'
$ws.Range("E49").Value = '  -5.65%  '
$ws.Range("D50").Value = '2.86'
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("D51").Value = '45.47'
$ws.Range("E51").Value = '  -2.69%  '
